$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current "Total" row (row 16), pushing Total to row 17
$ws.Rows.Item(16).Insert()

# Copy formatting from the row above (row 15, the last data row) into the new row 16
$ws.Range("A15:F15").Copy()
$ws.Range("A16:F16").PasteSpecial(-4122)

# Populate the new data row (2023-12-22 shift)
$ws.Range("A16").Value = 45282
$ws.Range("B16").Value = 0.625
$ws.Range("C16").Value = 0.916666666666667
$ws.Range("D16").Formula = "=(C16<B16)+C16-B16"
$ws.Range("E16").Value = 10
$ws.Range("F16").Formula = "=(D16*24)*E16"
$ws.Range("F16").NumberFormat = "General"

# Fix up the Total row's SUM ranges to include the new row
$ws.Range("D17").Formula = "=SUM(D2:D16)"
$ws.Range("F17").Formula = "=SUM(F2:F16)"

# Move the active selection to A17, matching the saved view state
$ws.Range("A17").Select()
